$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column A (shifts all existing data +2 cols)
$ws.Range("A1:B1").EntireColumn.Insert()

# 2. New header cells for the inserted columns
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Pais"

# 3. New data cells for row 2 in the inserted columns
$ws.Range("A2").Value = "19 octubre 2023"
$ws.Range("B2").Value = "México"

# 4. Replace the (shifted) company identity block in row 2
$ws.Range("C2").Value = "Swiss Steel Mexico SA de CV"
$ws.Range("D2").Value = "Swiss Steel Mexico"
$ws.Range("E2").Value = "`tSSB9512118M1"
$ws.Range("F2").Value = "Swiss Steel International  "

# Old Porcentaje1 / Accionistas2 / Porcentaje2 values (now at G2:I2) are gone entirely
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()
$ws.Range("I2").Clear()

# 5. Replace the (shifted) executives block in row 2
$ws.Range("P2").Value = "Thiery Jean Denis Cremailh "
$ws.Range("Q2").Value = "Presidente"
$ws.Range("R2").Value = "Sara Toriz Escamilla"
$ws.Range("S2").Value = "Secretario"
$ws.Range("T2").Value = "Fermin Huerta Rodriguez"
$ws.Range("U2").Value = "Comisario"
$ws.Range("V2").Value = "Jose Antonio Flores Muñoz"
$ws.Range("W2").Value = "Apoderado"

# 6. Replace the (shifted) Primarios block in row 2
$ws.Range("AF2").Value = "Swiss Steel International  "
$ws.Range("AG2").Value = "Thiery Jean Denis Cremailh "

# Old Primarios values that used to live at AH2:AI2 are gone entirely
$ws.Range("AH2").Clear()
$ws.Range("AI2").Clear()

# 7. Update the view: select B5, no frozen/scrolled top-left cell
$ws.Range("B5").Select()
